$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.735.11"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "1.895.06"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.01"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4940"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2966"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06822"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("D10").Value = "1.894.90"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.27"
$ws.Range("E11").Value = "  +3.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "92.54"
$ws.Range("E12").Value = "  +7.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07272"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.112"
$ws.Range("E14").Value = "  +4.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6803"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "30.696.70"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007992"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.33"
$ws.Range("E18").Value = "  +4.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "2.139.43"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.864"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "193.67"
$ws.Range("E23").Value = "  +36.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.089"
$ws.Range("E24").Value = "  +7.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.448"
$ws.Range("E25").Value = "  +4.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.81"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.26"
$ws.Range("E27").Value = "  +12.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.927"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.405"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.361"
$ws.Range("E30").Value = "  +3.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09022"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.045"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7474"
$ws.Range("E34").Value = "  +4.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.128"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.737"
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01867"
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.685"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.168"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9443"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4443"
$ws.Range("E41").Value = "  +4.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.50"
$ws.Range("E42").Value = "  +3.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.768"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.698"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("E46").Value = "  +6.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05865"
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.773"
$ws.Range("E48").Value = "  +5.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.438"
$ws.Range("E49").Value = "  +7.60%  "
$ws.Range("E50").Value = "  +4.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.65"
$ws.Range("E51").Value = "  +3.53%  "
